$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-29 Tuesday" "2025-04-30 Wednesday"
Replace-Text "769×2=" "511×8="
Replace-Text "154×5=" "104×4="
Replace-Text "208×5=" "738×2="
Replace-Text "422×8=" "172×2="
Replace-Text "765×7=" "404×2="
Replace-Text "222×4=" "479×8="
Replace-Text "559×9=" "335×5="
Replace-Text "269×6=" "101×9="
Replace-Text "807×8=" "287×8="
Replace-Text "400×5=" "839×6="
Replace-Text "413×2=" "998×9="
Replace-Text "407×8=" "940×7="
Replace-Text "251×2=" "710×8="
Replace-Text "440×3=" "410×9="
Replace-Text "727×4=" "540×9="
Replace-Text "321×6=" "300×5="
Replace-Text "224×7=" "113×5="
Replace-Text "393×4=" "700×8="
Replace-Text "678×2=" "267×5="
Replace-Text "181×4=" "433×9="
Replace-Text "349×8=" "847×7="
Replace-Text "125×9=" "631×8="
Replace-Text "194×3=" "357×2="
Replace-Text "320×7=" "795×3="
Replace-Text "980×5=" "709×6="
